$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.819.53'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '2.478.97'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'" + '316.08'
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').Value = "'" + '104.00'
$ws.Range('E6').Value = '  -4.96%  '
$ws.Range('E7').Value = '  -3.12%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -3.69%  '
$ws.Range('D10').Value = "'" + '38.62'
$ws.Range('E10').Value = '  -4.67%  '
$ws.Range('D11').Value = "'" + '20.40'
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').Value = "'" + '7.01'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').Value = '2.867.63'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '2.495.62'
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('D17').Value = "'" + '0.820'
$ws.Range('E17').Value = '  -4.01%  '
$ws.Range('D18').Value = '47.749.56'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').Value = "'" + '2.91'
$ws.Range('E19').Value = '  +7.24%  '
$ws.Range('D20').Value = "'" + '12.60'
$ws.Range('E20').Value = '  -5.37%  '
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').Value = '0.0₃0924'
$ws.Range('E22').Value = '  -2.39%  '
$ws.Range('D23').Value = "'" + '278.72'
$ws.Range('E23').Value = '  +5.34%  '
$ws.Range('D24').Value = "'" + '70.76'
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('E25').Value = '  -3.96%  '
$ws.Range('D26').Value = "'" + '0.999'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = "'" + '25.58'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').Value = "'" + '2.23'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').Value = "'" + '9.54'
$ws.Range('E29').Value = '  -5.73%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = "'" + '0.136'
$ws.Range('E30').Value = '  -5.04%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'" + '34.41'
$ws.Range('E31').Value = '  -4.30%  '
$ws.Range('D32').Value = "'" + '49.21'
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('E34').Value = '  -5.29%  '
$ws.Range('D35').Value = "'" + '5.22'
$ws.Range('E35').Value = '  -3.52%  '
$ws.Range('E36').Value = '  -3.33%  '
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('D38').Value = "'" + '4.46'
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('D39').Value = "'" + '2.85'
$ws.Range('E39').Value = '  -4.89%  '
$ws.Range('D40').Value = "'" + '122.47'
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Value = "'" + '21.78'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('D44').Value = "'" + '0.0297'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').Value = '1.990.52'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').Value = "'" + '3.11'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('E48').Value = '  -3.48%  '
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').Value = "'" + '5.08'
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('D51').Value = "'" + '78.84'
$ws.Range('E51').Value = '  -0.28%  '
